# Fix typo in header: "incricao_estadual" -> "inscricao_estadual"
# (also drops the cell's style/format, matching the saved state)
# and restore the last-used cell selection (N13) on the "Entidade" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entidade")

$ws.Range("E1").ClearFormats()
$ws.Range("E1").Value = "inscricao_estadual"

$ws.Activate()
$ws.Range("N13").Select()
